# Update "想去人数" (number of people interested) counts per the latest
# scrape output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1585
$wsExhibit.Range("F5").Value = 729
$wsExhibit.Range("F6").Value = 37

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 12

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1585
$wsAll.Range("F5").Value = 12
$wsAll.Range("F6").Value = 729
$wsAll.Range("F7").Value = 37
